$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells that look like plain numbers to stay text,
# matching the source data (coinranking price strings are text, not numeric).

$ws.Range("D2").Value = "42.483.56"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.517.83"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.68"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.57"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  -1.08%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.06"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.20"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("D14").Value = "2.905.17"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "2.535.37"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.04"
$ws.Range("E16").Value = "  -6.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.811"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "42.467.03"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.55"
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.03"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.35"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.55"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.85"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  -3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.38"
$ws.Range("E28").Value = "  -4.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.51"
$ws.Range("E30").Value = "  -5.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.83"
$ws.Range("E31").Value = "  +2.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.40"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.63"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0779"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.12"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -3.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.35"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.116"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.10"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "2.026.70"
$ws.Range("E44").Value = "  +4.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0293"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.20"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "2.763.96"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "79.22"
$ws.Range("E49").Value = "  -1.51%  "

# Row 50/51 swap: Algorand <-> ordi
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.11"
$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.186"
$ws.Range("E51").Value = "  -2.58%  "
